$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.09494813238296555
$ws.Range("C2").Value = 0.5332219804412588
$ws.Range("D2").Value = 0.4006434568944247
$ws.Range("E2").Value = 0.6329640249606803
$ws.Range("F2").Value = 0.6398668578903548
$ws.Range("G2").Value = 23

$ws.Range("B3").Value = 0.01602953393462571
$ws.Range("C3").Value = 0.5143741298588096
$ws.Range("D3").Value = 0.4182794719077274
$ws.Range("E3").Value = 0.6467452913688103
$ws.Range("F3").Value = 0.6617615586588519
$ws.Range("G3").Value = 22

$ws.Range("B4").Value = 0.1085716555968761
$ws.Range("C4").Value = 0.452885312363856
$ws.Range("D4").Value = 0.3549274518830394
$ws.Range("E4").Value = 0.5957578802525733
$ws.Range("F4").Value = 0.6002471406497429
$ws.Range("G4").Value = 21

$ws.Range("B5").Value = 0.1299845358620406
$ws.Range("C5").Value = 0.608575057432829
$ws.Range("D5").Value = 0.5363125488530067
$ws.Range("E5").Value = 0.7323336322012028
$ws.Range("F5").Value = 0.7394283491068494
$ws.Range("G5").Value = 20

$ws.Range("B6").Value = 0.07656759945532947
$ws.Range("C6").Value = 0.3700814829430278
$ws.Range("D6").Value = 0.2616870517883947
$ws.Range("E6").Value = 0.5115535668807273
$ws.Range("F6").Value = 0.5196507713807427
$ws.Range("G6").Value = 19

$ws.Range("B7").Value = 0.2093201528469608
$ws.Range("C7").Value = 0.5868146440123966
$ws.Range("D7").Value = 0.5026012083360064
$ws.Range("E7").Value = 0.7089437271998437
$ws.Range("F7").Value = 0.6969746841155865
$ws.Range("G7").Value = 18

$ws.Range("B8").Value = 0.1086669261956075
$ws.Range("C8").Value = 0.4904301031709458
$ws.Range("D8").Value = 0.3878951139765857
$ws.Range("E8").Value = 0.6228122622240073
$ws.Range("F8").Value = 0.6321329183393875
$ws.Range("G8").Value = 17
